$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '64.596.88'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.84%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.158.73'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +3.61%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '565.44'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.96%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '143.17'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +4.28%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.149.95'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +3.59%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.497'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +2.64%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.79'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +6.07%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.155'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +2.85%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.468'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +3.36%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '36.97'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +4.78%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000223'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +2.79%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.656.85'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +3.33%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.611.86'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.73%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.149.45'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +3.10%  '
$ws.Range('E18').Value = '  +1.39%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '518.71'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +7.60%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.88'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +5.21%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.11'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +4.43%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.720'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +6.03%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.50'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +5.72%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '12.83'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +4.93%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '79.25'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.30%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.84'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +15.87%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.84'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +5.70%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.16'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +4.91%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.998'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.36%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '26.71'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +4.15%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.61'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.57%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.14'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +2.94%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '555.06'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -4.42%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '6.11'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +4.20%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.41'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.36%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '53.98'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +3.58%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0436'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +8.40%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0828'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +5.56%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.171.24'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +8.80%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.123'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +5.02%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.76'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.17%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.34'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +2.38%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.266'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +10.95%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.23'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +8.49%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '25.41'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +3.81%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '121.05'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +2.95%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.109'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.69%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0₃0520'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.11'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +4.04%  '
